$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.147.70"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "'3.372.59"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D5").Value = "'405.88"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'134.93"
$ws.Range("E6").Value = "  +9.79%  "
$ws.Range("D7").Value = "'0.592"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.674"
$ws.Range("E9").Value = "  +5.35%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'42.82"
$ws.Range("E11").Value = "  +4.13%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "'3.886.35"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "'8.35"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'19.60"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "'3.381.44"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'61.019.60"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "'1.02"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'11.01"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("D21").Value = "'3.22"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").Value = "'83.62"
$ws.Range("E22").Value = "  +9.26%  "
$ws.Range("D23").Value = "'312.87"
$ws.Range("E23").Value = "  +4.70%  "
$ws.Range("D24").Value = "'12.74"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  +11.47%  "
$ws.Range("D27").Value = "'8.41"
$ws.Range("E27").Value = "  +10.21%  "
$ws.Range("D28").Value = "'29.50"
$ws.Range("E28").Value = "  -3.93%  "
$ws.Range("D29").Value = "'7.46"
$ws.Range("E29").Value = "  -7.38%  "
$ws.Range("D30").Value = "'0.173"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'11.31"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'41.26"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "'0.0481"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "'52.15"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("D40").Value = "'2.92"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "'137.35"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "'0.291"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").Value = "'16.63"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").Value = "'21.44"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'2.123.56"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -5.09%  "
$ws.Range("E51").Value = "  -0.39%  "
